$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) to the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 and IF data columns for rows 2-83
$data = @(
    @(2, 4, 6),
    @(3, 6, 7),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 5, 5),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 6, 7),
    @(12, 1, 2),
    @(13, 8, 8),
    @(14, 1, 2),
    @(15, 6, 7),
    @(16, 6, 7),
    @(17, 6, 7),
    @(18, 7, 7),
    @(19, 8, 9),
    @(20, 8, 9),
    @(21, 5, 6),
    @(22, 7, 7),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 6, 6),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 6, 7),
    @(38, 8, 8),
    @(39, 7, 8),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 6, 7),
    @(44, 7, 8),
    @(45, 8, 8),
    @(46, 7, 8),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 8, 8),
    @(51, 8, 8),
    @(52, 6, 7),
    @(53, 6, 7),
    @(54, 7, 8),
    @(55, 9, 9),
    @(56, 8, 8),
    @(57, 8, 8),
    @(58, 7, 7),
    @(59, 7, 8),
    @(60, 8, 8),
    @(61, 8, 8),
    @(62, 8, 8),
    @(63, 8, 8),
    @(64, 8, 9),
    @(65, 8, 8),
    @(66, 7, 8),
    @(67, 7, 7),
    @(68, 7, 8),
    @(69, 8, 8),
    @(70, 8, 8),
    @(71, 9, 9),
    @(72, 7, 8),
    @(73, 8, 8),
    @(74, 6, 6),
    @(75, 7, 7),
    @(76, 8, 8),
    @(77, 6, 6),
    @(78, 8, 8),
    @(79, 7, 7),
    @(80, 8, 8),
    @(81, 5, 5),
    @(82, 4, 4),
    @(83, 4, 4)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
